$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("station")
$ws.Activate()

# The "Stop time (min)" column header is simplified to "Stop time" now that
# the sheet also reports electrified driving time, so the old "(min)"
# qualifier is no longer accurate/needed.
$ws.Range("D1").Value = "Stop time"

# Re-fit the header row now that the (shorter) text no longer wraps to two
# lines, and put the selection back on the cell we just edited.
$ws.Rows.Item(1).EntireRow.AutoFit()
$ws.Range("D1").Select()
